# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# Values are written as literal text (not re-parsed as numbers/dates) so strings
# like "596.29" or "1.00" keep their exact display text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"   # force text interpretation before assignment
    $r.Value = $val
    $r.Style = "Normal"    # drop the temporary text format again
}

Set-TextValue "D2" "63.824.38"
Set-TextValue "E2" "  -0.08%  "
Set-TextValue "D3" "2.625.84"
Set-TextValue "E3" "  +0.02%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "596.29"
Set-TextValue "E5" "  -0.04%  "
Set-TextValue "D6" "151.40"
Set-TextValue "E6" "  +0.83%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "0.587"
Set-TextValue "E8" "  -0.17%  "
Set-TextValue "E9" "  +4.76%  "
Set-TextValue "D10" "5.84"
Set-TextValue "E10" "  +3.04%  "
Set-TextValue "D11" "0.397"
Set-TextValue "E11" "  +3.83%  "
Set-TextValue "E12" "  +0.81%  "
Set-TextValue "D13" "27.95"
Set-TextValue "E13" "  +0.91%  "
Set-TextValue "D14" "3.098.34"
Set-TextValue "E14" "  +0.04%  "
Set-TextValue "D15" "63.735.88"
Set-TextValue "D16" "0.0000168"
Set-TextValue "E16" "  +13.22%  "
Set-TextValue "D17" "2.603.37"
Set-TextValue "E17" "  -1.34%  "
Set-TextValue "D18" "12.18"
Set-TextValue "E18" "  -1.41%  "
Set-TextValue "D19" "4.80"
Set-TextValue "E19" "  +3.60%  "
Set-TextValue "D20" "347.58"
Set-TextValue "E20" "  -0.60%  "
Set-TextValue "D21" "7.01"
Set-TextValue "E21" "  +1.50%  "
Set-TextValue "E22" "  +0.25%  "
Set-TextValue "D23" "67.48"
Set-TextValue "E23" "  +1.77%  "
Set-TextValue "D24" "1.69"
Set-TextValue "E24" "  -2.89%  "
Set-TextValue "E25" "  +0.68%  "
Set-TextValue "D26" "9.13"
Set-TextValue "E26" "  -0.87%  "
Set-TextValue "D27" "8.33"
Set-TextValue "E27" "  +1.69%  "
Set-TextValue "D28" "554.58"
Set-TextValue "E28" "  -1.33%  "
Set-TextValue "D29" "0.163"
Set-TextValue "E29" "  -1.23%  "
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "B31" "PEPE"
Set-TextValue "C31" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D31" "0.0₃0903"
Set-TextValue "E31" "  +6.78%  "
Set-TextValue "B32" "PancakeSwap"
Set-TextValue "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "2.06"
Set-TextValue "E32" "  +1.13%  "
Set-TextValue "E33" "  +4.68%  "
Set-TextValue "D34" "5.38"
Set-TextValue "E34" "  +3.25%  "
Set-TextValue "D35" "6.12"
Set-TextValue "E35" "  +0.81%  "
Set-TextValue "B36" "PolygonEcosystemToken"
Set-TextValue "C36" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D36" "0.418"
Set-TextValue "E36" "  +2.26%  "
Set-TextValue "B37" "Monero"
Set-TextValue "C37" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "164.38"
Set-TextValue "E37" "  -3.08%  "
Set-TextValue "B38" "EthereumClassic"
Set-TextValue "C38" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D38" "19.95"
Set-TextValue "E38" "  +2.90%  "
Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.99"
Set-TextValue "E39" "  +1.95%  "
Set-TextValue "E40" "  -0.06%  "
Set-TextValue "E41" "  -0.06%  "
Set-TextValue "D42" "168.56"
Set-TextValue "E42" "  -1.18%  "
Set-TextValue "D43" "4.08"
Set-TextValue "E43" "  +3.89%  "
Set-TextValue "D44" "23.47"
Set-TextValue "E44" "  +9.45%  "
Set-TextValue "D45" "0.0586"
Set-TextValue "E45" "  -2.07%  "
Set-TextValue "D46" "2.19"
Set-TextValue "E46" "  +10.86%  "
Set-TextValue "D47" "0.638"
Set-TextValue "E47" "  +1.18%  "
Set-TextValue "E48" "  +2.32%  "
Set-TextValue "D49" "0.0970"
Set-TextValue "E49" "  +0.20%  "
Set-TextValue "D50" "19.28"
Set-TextValue "E50" "  +0.26%  "
Set-TextValue "D51" "0.0₆0233"
Set-TextValue "E51" "  +19.42%  "
